# The commit inserts a brand-new price-report row right before the
# existing row 247 ("Fruta / hortaliza, semanal" — a new week's record
# for Zanahoria @ Macroferia Regional de Talca), pushing every
# following row down by one (old row 352 becomes new row 353).
#
# Inserting a whole row (rather than just writing into row 247) is what
# reproduces the diff: every row from 247 to 352 keeps its original
# values but slides down one position, dimension grows from R352 to
# R353, and the newly opened row 247 gets the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 247..352 down to 248..353, opening up a blank row 247.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new record.
$ws.Range("A247").Value = 5
$ws.Range("B247").Value = "Macroferia Regional de Talca"
$ws.Range("C247").Value = "Maule"
$ws.Range("D247").Value = 44726
$ws.Range("E247").Value = 7
$ws.Range("F247").Value = 100114013
$ws.Range("G247").Value = "Zanahoria"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 500
$ws.Range("K247").Value = 6500
$ws.Range("L247").Value = 6500
$ws.Range("M247").Value = 6500
$ws.Range("N247").Value = "$/saco 20 kilos"
$ws.Range("O247").Value = "Región de Ñuble"
$ws.Range("P247").Value = 325
$ws.Range("Q247").Value = 20
$ws.Range("R247").Value = "Hortaliza"
